# Update cryptos list figures (price/volume) per Sun Oct 27 07:44:17 UTC 2024 GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.182.39'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.21%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.479.77'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.54%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.61%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.18'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.11%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  +0.14%  '
$ws.Range("E9").Value = '  +2.53%  '
$ws.Range("E10").Value = '  +1.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '4.94'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.43%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.334'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.30%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.45'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '67.054.89'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.32%  '
$ws.Range("E16").Value = '  +0.54%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.494.64'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.36%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.60'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.83%  '
$ws.Range("E19").Value = '  -2.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '350.55'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.12%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.03'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.46%  '
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("E23").Value = '  -0.29%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.24'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.22%  '
$ws.Range("E25").Value = '  +2.15%  '
$ws.Range("E26").Value = '  +1.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.607.76'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.51%  '
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("E29").Value = '  +0.59%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '503.91'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.07%  '
$ws.Range("E31").Value = '  -0.21%  '
$ws.Range("E32").Value = '  +0.41%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.76'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.31%  '
$ws.Range("E34").Value = '  +0.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '162.06'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.43%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.119'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.47%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.69'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.41%  '
$ws.Range("E38").Value = '  -1.07%  '
$ws.Range("E39").Value = '  -1.35%  '
$ws.Range("E40").Value = '  +0.04%  '
$ws.Range("E41").Value = '  +1.68%  '
$ws.Range("E42").Value = '  +0.50%  '
$ws.Range("E44").Value = '  +1.84%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '142.94'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.59%  '
$ws.Range("B46").Value = 'Filecoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.48'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.64%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₆0260'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.70%  '
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("E49").Value = '  +0.59%  '
$ws.Range("E50").Value = '  -0.36%  '
$ws.Range("E51").Value = '  +0.59%  '
